$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Tue Nov 12 18:03:53 EST 2024"
$ws.Range("B3").Value = "Tue Nov 12 18:04:06 EST 2024"
$ws.Range("B4").Value = "Tue Nov 12 18:04:18 EST 2024"
$ws.Range("B5").Value = "Tue Nov 12 18:04:31 EST 2024"
